# Auto-generated Excel COM-interop script applying numeric updates
# to the Chocobo_Profits workbook sheets, per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2985
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 2985
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 2985
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -3335
$ws.Range("H93").Value = 22755.172
$ws.Range("J93").Value = 22755.172
$ws.Range("L93").Value = 22755.172
$ws.Range("N93").Value = -27747.172
$ws.Range("H99").Value = 1251.5555
$ws.Range("I99").Value = 594
$ws.Range("J99").Value = 2566.6667
$ws.Range("K99").Value = 1782
$ws.Range("L99").Value = 7700.000100000001
$ws.Range("M99").Value = -284
$ws.Range("N99").Value = -10696.0001
$ws.Range("H101").Value = 3982.7778
$ws.Range("J101").Value = 4459.625
$ws.Range("L101").Value = 13378.875
$ws.Range("N101").Value = -16622.875
$ws.Range("H105").Value = 30032.834
$ws.Range("J105").Value = 30032.834
$ws.Range("L105").Value = 30032.834
$ws.Range("N105").Value = -37020.834
$ws.Range("H113").Value = 11468.333
$ws.Range("I113").Value = 2405
$ws.Range("J113").Value = 16000
$ws.Range("K113").Value = 2405
$ws.Range("L113").Value = 16000
$ws.Range("M113").Value = 849
$ws.Range("N113").Value = -22508
$ws.Range("H115").Value = 1554
$ws.Range("I115").Value = 1554
$ws.Range("K115").Value = 4662
$ws.Range("M115").Value = -3095
$ws.Range("H116").Value = 507520
$ws.Range("I116").Value = 1002899.4
$ws.Range("K116").Value = 1002899.4
$ws.Range("M116").Value = -999457.4
$ws.Range("H118").Value = 893.3333
$ws.Range("I118").Value = 340
$ws.Range("J118").Value = 2000
$ws.Range("K118").Value = 1020
$ws.Range("L118").Value = 6000
$ws.Range("M118").Value = 637
$ws.Range("N118").Value = -9314
$ws.Range("H125").Value = 2535.6667
$ws.Range("J125").Value = 2535.6667
$ws.Range("L125").Value = 22821.0003
$ws.Range("N125").Value = -27741.0003
$ws.Range("H127").Value = 1865.875
$ws.Range("I127").Value = 543.5
$ws.Range("J127").Value = 2306.6667
$ws.Range("K127").Value = 1630.5
$ws.Range("L127").Value = 6920.000100000001
$ws.Range("M127").Value = 3329.5
$ws.Range("N127").Value = -16840.0001
$ws.Range("H129").Value = 836.49
$ws.Range("J129").Value = 864.7578999999999
$ws.Range("L129").Value = 2594.2737
$ws.Range("N129").Value = -12594.2737
$ws.Range("H132").Value = 45462060
$ws.Range("I132").Value = 58831770
$ws.Range("J132").Value = 5040
$ws.Range("K132").Value = 176495310
$ws.Range("L132").Value = 15120
$ws.Range("M132").Value = -176492780
$ws.Range("N132").Value = -20180
$ws.Range("H138").Value = 2657.48
$ws.Range("I138").Value = 854.8333
$ws.Range("J138").Value = 2903.2954
$ws.Range("K138").Value = 2564.4999
$ws.Range("L138").Value = 8709.886200000001
$ws.Range("M138").Value = 2575.5001
$ws.Range("N138").Value = -18989.8862
$ws.Range("H141").Value = 47996.547
$ws.Range("I141").Value = 68749.60000000001
$ws.Range("J141").Value = 3525.7144
$ws.Range("K141").Value = 206248.8
$ws.Range("L141").Value = 10577.1432
$ws.Range("M141").Value = -201068.8
$ws.Range("N141").Value = -20937.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 893.3333
$ws.Range("I2").Value = 850
$ws.Range("J2").Value = 980
$ws.Range("K2").Value = 850
$ws.Range("L2").Value = 980
$ws.Range("M2").Value = -737
$ws.Range("N2").Value = -1206
$ws.Range("H61").Value = 1623
$ws.Range("I61").Value = 1525.1333
$ws.Range("K61").Value = 1525.1333
$ws.Range("M61").Value = -1313.1333
$ws.Range("H93").Value = 24500
$ws.Range("J93").Value = 24500
$ws.Range("L93").Value = 24500
$ws.Range("N93").Value = -29492
$ws.Range("H116").Value = 893.3333
$ws.Range("I116").Value = 850
$ws.Range("J116").Value = 980
$ws.Range("K116").Value = 850
$ws.Range("L116").Value = 980
$ws.Range("M116").Value = 1444
$ws.Range("N116").Value = -5568
$ws.Range("H122").Value = 2618.25
$ws.Range("I122").Value = 849.4286
$ws.Range("K122").Value = 2548.2858
$ws.Range("M122").Value = -98.28579999999965
$ws.Range("H132").Value = 4242
$ws.Range("I132").Value = 1266.6666
$ws.Range("J132").Value = 5729.6665
$ws.Range("K132").Value = 3799.9998
$ws.Range("L132").Value = 17188.9995
$ws.Range("M132").Value = -1269.9998
$ws.Range("N132").Value = -22248.9995
$ws.Range("H136").Value = 1623
$ws.Range("I136").Value = 1525.1333
$ws.Range("K136").Value = 4575.3999
$ws.Range("M136").Value = -2025.3999
$ws.Range("H137").Value = 42255
$ws.Range("J137").Value = 42255
$ws.Range("L137").Value = 42255
$ws.Range("N137").Value = -52455

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 893.3333
$ws.Range("I3").Value = 850
$ws.Range("J3").Value = 980
$ws.Range("K3").Value = 850
$ws.Range("L3").Value = 980
$ws.Range("M3").Value = -736
$ws.Range("N3").Value = -1208
$ws.Range("H134").Value = 2141.7
$ws.Range("I134").Value = 1402.0588
$ws.Range("J134").Value = 6333
$ws.Range("K134").Value = 4206.1764
$ws.Range("L134").Value = 18999
$ws.Range("M134").Value = -1671.1764
$ws.Range("N134").Value = -24069
$ws.Range("H137").Value = 38662.5
$ws.Range("J137").Value = 38662.5
$ws.Range("L137").Value = 38662.5
$ws.Range("N137").Value = -48862.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6664.5
$ws.Range("I31").Value = 1208.1765
$ws.Range("J31").Value = 13799.692
$ws.Range("K31").Value = 1208.1765
$ws.Range("L31").Value = 13799.692
$ws.Range("M31").Value = -913.1765
$ws.Range("N31").Value = -14389.692
$ws.Range("H34").Value = 6664.5
$ws.Range("I34").Value = 1208.1765
$ws.Range("J34").Value = 13799.692
$ws.Range("K34").Value = 1208.1765
$ws.Range("L34").Value = 13799.692
$ws.Range("M34").Value = -1006.1765
$ws.Range("N34").Value = -14203.692
$ws.Range("H120").Value = 39999
$ws.Range("J120").Value = 39999
$ws.Range("L120").Value = 39999
$ws.Range("N120").Value = -47257
$ws.Range("H137").Value = 33355.43
$ws.Range("J137").Value = 33355.43
$ws.Range("L137").Value = 33355.43
$ws.Range("N137").Value = -43555.43
$ws.Range("H139").Value = 38900
$ws.Range("J139").Value = 38900
$ws.Range("L139").Value = 38900
$ws.Range("N139").Value = -49180

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 470
$ws.Range("I7").Value = 440
$ws.Range("J7").Value = 500
$ws.Range("K7").Value = 1320
$ws.Range("L7").Value = 1500
$ws.Range("M7").Value = -1208
$ws.Range("N7").Value = -1724
$ws.Range("H12").Value = 72.48148
$ws.Range("I12").Value = 18.333334
$ws.Range("J12").Value = 99.55556
$ws.Range("K12").Value = 55.000002
$ws.Range("L12").Value = 298.66668
$ws.Range("M12").Value = 117.999998
$ws.Range("N12").Value = -644.66668
$ws.Range("H43").Value = 9000
$ws.Range("J43").Value = 9000
$ws.Range("L43").Value = 27000
$ws.Range("N43").Value = -27228
$ws.Range("H80").Value = 9166.333000000001
$ws.Range("J80").Value = 8749.5
$ws.Range("L80").Value = 26248.5
$ws.Range("N80").Value = -28120.5
$ws.Range("H83").Value = 9166.333000000001
$ws.Range("J83").Value = 8749.5
$ws.Range("L83").Value = 78745.5
$ws.Range("N83").Value = -88105.5
$ws.Range("H92").Value = 763.5714
$ws.Range("I92").Value = 591.53845
$ws.Range("J92").Value = 3000
$ws.Range("K92").Value = 1774.61535
$ws.Range("L92").Value = 9000
$ws.Range("M92").Value = -526.61535
$ws.Range("N92").Value = -11496
$ws.Range("H113").Value = 3677194.5
$ws.Range("I113").Value = 594.4286
$ws.Range("K113").Value = 1783.2858
$ws.Range("M113").Value = 386.7142000000001
$ws.Range("H129").Value = 1984.125
$ws.Range("J129").Value = 2429
$ws.Range("L129").Value = 7287
$ws.Range("N129").Value = -17287

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 23922.223
$ws.Range("J46").Value = 23922.223
$ws.Range("L46").Value = 23922.223
$ws.Range("N46").Value = -24234.223
$ws.Range("H102").Value = 2265.1516
$ws.Range("I102").Value = 1701.9615
$ws.Range("K102").Value = 1701.9615
$ws.Range("M102").Value = -79.96149999999989
$ws.Range("H132").Value = 5172.6665
$ws.Range("I132").Value = 4489.8667
$ws.Range("J132").Value = 6879.6665
$ws.Range("K132").Value = 13469.6001
$ws.Range("L132").Value = 20638.9995
$ws.Range("M132").Value = -10939.6001
$ws.Range("N132").Value = -25698.9995
$ws.Range("H137").Value = 40220
$ws.Range("J137").Value = 40220
$ws.Range("L137").Value = 40220
$ws.Range("N137").Value = -50420

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5913
$ws.Range("I40").Value = 4614.857
$ws.Range("K40").Value = 4614.857
$ws.Range("M40").Value = -4478.857

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2970.5264
$ws.Range("I81").Value = 2539.4285
$ws.Range("J81").Value = 4177.6
$ws.Range("K81").Value = 5078.857
$ws.Range("L81").Value = 8355.200000000001
$ws.Range("M81").Value = -4017.857
$ws.Range("N81").Value = -10477.2
$ws.Range("H84").Value = 2970.5264
$ws.Range("I84").Value = 2539.4285
$ws.Range("J84").Value = 4177.6
$ws.Range("K84").Value = 25394.285
$ws.Range("L84").Value = 41776
$ws.Range("M84").Value = -20090.285
$ws.Range("N84").Value = -52384
$ws.Range("H107").Value = 923.75
$ws.Range("I107").Value = 851
$ws.Range("J107").Value = 996.5
$ws.Range("K107").Value = 2553
$ws.Range("L107").Value = 2989.5
$ws.Range("M107").Value = -633
$ws.Range("N107").Value = -6829.5
$ws.Range("H122").Value = 12333.333
$ws.Range("I122").Value = 11000
$ws.Range("K122").Value = 33000
$ws.Range("M122").Value = -30550
$ws.Range("H132").Value = 15158382
$ws.Range("I132").Value = 8754.846
$ws.Range("J132").Value = 37041176
$ws.Range("K132").Value = 26264.538
$ws.Range("L132").Value = 111123528
$ws.Range("M132").Value = -23734.538
$ws.Range("N132").Value = -111128588
